$d = $word.ActiveDocument

# Merge the five runs of the title paragraph ("Testing" / " " / "custom" / " " / "properties")
# into a single run containing "Testing custom properties".
$d.Content.Find.Execute("Testing custom properties", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Testing custom properties", 2)

# Merge the three runs of the author paragraph ("A." / " " / "M.") into a single run
# containing "A. M.".
$d.Content.Find.Execute("A. M.", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "A. M.", 2)
